$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (46081 -> 46082, i.e. 2026-02-28 -> 2026-03-01) for every data row
# (rows 2 through 496).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 496 }

$ws.Range("C2:C$lastRow").Value = 46082
